# The deck's custom "Integral" theme (ppt/theme/theme1.xml, used by the
# slide master / all slides) is switched to the default PowerPoint
# "Office Theme" colour palette (the palette that used to live only in
# ppt/theme/theme2.xml, the notes-master theme).
#
# The DrawingML <a:clrScheme> defines 12 colours, in this fixed order:
#   1 dk1  2 lt1  3 dk2  4 lt2  5 accent1  6 accent2  7 accent3
#   8 accent4  9 accent5  10 accent6  11 hlink  12 folHlink
#
# PowerPoint's object model exposes that full 12-slot palette via
# Slide.ThemeColorScheme (Slide/SlideRange/CustomLayout all resolve to
# the same underlying presentation theme part), where each ThemeColor's
# .RGB is a standard Win32 COLORREF (0x00BBGGRR, i.e. byte-reversed
# from the "RRGGBB" hex most people think in).

function ColorRefFromHex($hex) {
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    return ($b * 65536) + ($g * 256) + $r
}

# Target palette = the standard Office theme colours (RRGGBB), in
# clrScheme order.
$officeThemeHex = @(
    "000000", # dk1
    "FFFFFF", # lt1
    "44546A", # dk2
    "E7E6E6", # lt2
    "5B9BD5", # accent1
    "ED7D31", # accent2
    "A5A5A5", # accent3
    "FFC000", # accent4
    "4472C4", # accent5
    "70AD47", # accent6
    "0563C1", # hlink
    "954F72"  # folHlink
)

$p = $ppt.ActivePresentation

# Any slide's ThemeColorScheme reaches into the shared presentation
# theme part (ppt/theme/theme1.xml) used by the slide master/layouts.
$themeColors = $p.Slides.Item(1).ThemeColorScheme

for ($i = 1; $i -le $themeColors.Count; $i++) {
    $themeColors.Colors($i).RGB = ColorRefFromHex $officeThemeHex[$i - 1]
}
